$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-08-10 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-11 Monday", 2)

$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}

Set-CellText 1 1 "52÷2=26, 0"
Set-CellText 1 2 "54÷5=10, 4"
Set-CellText 1 3 "25÷3=8, 1"
Set-CellText 1 4 "22÷4=5, 2"
Set-CellText 1 5 "40÷9=4, 4"

Set-CellText 5 1 "26÷8=3, 2"
Set-CellText 5 2 "41÷5=8, 1"
Set-CellText 5 3 "40÷5=8, 0"
Set-CellText 5 4 "71÷9=7, 8"
Set-CellText 5 5 "27÷4=6, 3"

Set-CellText 9 1 "36÷7=5, 1"
Set-CellText 9 2 "16÷2=8, 0"
Set-CellText 9 3 "98÷3=32, 2"
Set-CellText 9 4 "78÷7=11, 1"
Set-CellText 9 5 "29÷5=5, 4"

Set-CellText 13 1 "74÷9=8, 2"
Set-CellText 13 2 "16÷8=2, 0"
Set-CellText 13 3 "26÷5=5, 1"
Set-CellText 13 4 "11÷6=1, 5"
Set-CellText 13 5 "67÷6=11, 1"

Set-CellText 17 1 "21÷4=5, 1"
Set-CellText 17 2 "41÷4=10, 1"
Set-CellText 17 3 "11÷7=1, 4"
Set-CellText 17 4 "26÷9=2, 8"
Set-CellText 17 5 "14÷3=4, 2"
